# Automatic update of files.
# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45180 (2023-09-11) to 45181 (2023-09-12).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row  # xlUp = -4162
if ($lastRow -lt 2) { $lastRow = 20 }

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45180) {
        $cell.Value2 = 45181
    }
}
